{"js": "// Replace the two-digit multiplication equations in the document's table\n// with the new values, per the commit's OOXML diff. Every original equation\n// string is unique in the document, so a direct search/replace keyed on the\n// exact old text is safe and will not clobber unrelated cells.\nconst replacements = [\n  [\"42\u00d777=3234\", \"94\u00d763=5922\"],\n  [\"11\u00d734=374\", \"11\u00d745=495\"],\n  [\"43\u00d739=1677\", \"17\u00d771=1207\"],\n  [\"98\u00d749=4802\", \"42\u00d773=3066\"],\n  [\"74\u00d738=2812\", \"50\u00d769=3450\"],\n  [\"67\u00d765=4355\", \"45\u00d751=2295\"],\n  [\"95\u00d748=4560\", \"81\u00d790=7290\"],\n  [\"72\u00d788=6336\", \"70\u00d715=1050\"],\n  [\"79\u00d780=6320\", \"89\u00d733=2937\"],\n  [\"34\u00d757=1938\", \"35\u00d789=3115\"],\n  [\"44\u00d711=484\", \"76\u00d769=5244\"],\n  [\"60\u00d773=4380\", \"23\u00d741=943\"],\n  [\"80\u00d758=4640\", \"20\u00d711=220\"],\n  [\"63\u00d774=4662\", \"76\u00d751=3876\"],\n  [\"16\u00d746=736\", \"55\u00d774=4070\"],\n  [\"57\u00d739=2223\", \"89\u00d740=3560\"],\n  [\"74\u00d723=1702\", \"59\u00d786=5074\"],\n  [\"69\u00d761=4209\", \"80\u00d719=1520\"],\n  [\"72\u00d778=5616\", \"33\u00d795=3135\"],\n  [\"94\u00d712=1128\", \"78\u00d754=4212\"],\n  [\"53\u00d738=2014\", \"60\u00d728=1680\"],\n  [\"81\u00d768=5508\", \"46\u00d733=1518\"],\n  [\"43\u00d770=3010\", \"45\u00d729=1305\"],\n  [\"95\u00d794=8930\", \"70\u00d724=1680\"],\n  [\"13\u00d769=897\", \"14\u00d758=812\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text to replace: ${oldText}`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the two-digit multiplication equations in the document's table\n# with the new values, per the commit's OOXML diff. Every original equation\n# string is unique in the document, so a direct Find/Replace keyed on the\n# exact old text is safe and will not clobber unrelated cells.\n$pairs = @(\n  @(\"42\u00d777=3234\", \"94\u00d763=5922\"),\n  @(\"11\u00d734=374\", \"11\u00d745=495\"),\n  @(\"43\u00d739=1677\", \"17\u00d771=1207\"),\n  @(\"98\u00d749=4802\", \"42\u00d773=3066\"),\n  @(\"74\u00d738=2812\", \"50\u00d769=3450\"),\n  @(\"67\u00d765=4355\", \"45\u00d751=2295\"),\n  @(\"95\u00d748=4560\", \"81\u00d790=7290\"),\n  @(\"72\u00d788=6336\", \"70\u00d715=1050\"),\n  @(\"79\u00d780=6320\", \"89\u00d733=2937\"),\n  @(\"34\u00d757=1938\", \"35\u00d789=3115\"),\n  @(\"44\u00d711=484\", \"76\u00d769=5244\"),\n  @(\"60\u00d773=4380\", \"23\u00d741=943\"),\n  @(\"80\u00d758=4640\", \"20\u00d711=220\"),\n  @(\"63\u00d774=4662\", \"76\u00d751=3876\"),\n  @(\"16\u00d746=736\", \"55\u00d774=4070\"),\n  @(\"57\u00d739=2223\", \"89\u00d740=3560\"),\n  @(\"74\u00d723=1702\", \"59\u00d786=5074\"),\n  @(\"69\u00d761=4209\", \"80\u00d719=1520\"),\n  @(\"72\u00d778=5616\", \"33\u00d795=3135\"),\n  @(\"94\u00d712=1128\", \"78\u00d754=4212\"),\n  @(\"53\u00d738=2014\", \"60\u00d728=1680\"),\n  @(\"81\u00d768=5508\", \"46\u00d733=1518\"),\n  @(\"43\u00d770=3010\", \"45\u00d729=1305\"),\n  @(\"95\u00d794=8930\", \"70\u00d724=1680\"),\n  @(\"13\u00d769=897\", \"14\u00d758=812\")\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
